$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "48.177.49"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "2.528.28"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "323.69"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "109.22"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.564"
$ws.Range("E9").Value = "  +5.02%  "
$ws.Range("D10").Value = "40.44"
$ws.Range("E10").Value = "  +3.40%  "
$ws.Range("E11").Value = "  +10.17%  "
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").Value = "2.923.87"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").Value = "2.527.48"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").Value = "0.861"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").Value = "48.155.16"
$ws.Range("D19").Value = "13.25"
$ws.Range("E19").Value = "  +3.08%  "
$ws.Range("D20").Value = "6.63"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").Value = "0.0₃0948"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("D22").Value = "2.72"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "72.49"
$ws.Range("E23").Value = "  +2.79%  "
$ws.Range("D24").Value = "270.25"
$ws.Range("E24").Value = "  +9.21%  "
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").Value = "26.26"
$ws.Range("E26").Value = "  +0.78%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").Value = "10.16"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("D29").Value = "0.147"
$ws.Range("E29").Value = "  +6.70%  "
$ws.Range("D30").Value = "35.31"
$ws.Range("E31").Value = "  -8.77%  "
$ws.Range("D32").Value = "49.79"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").Value = "20.00"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "0.0792"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "22.64"
$ws.Range("E41").Value = "  +6.60%  "
$ws.Range("D42").Value = "2.20"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").Value = "118.01"
$ws.Range("E43").Value = "  -2.75%  "
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "2.010.30"
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("E46").Value = "  +2.60%  "
$ws.Range("E47").Value = "  +6.34%  "
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("D51").Value = "80.29"
$ws.Range("E51").Value = "  +2.90%  "

